$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-6, column A/B are stored as *text* (not numbers) in the target
# workbook, even though several of the values look numeric. Prefixing
# with a leading apostrophe is the normal Excel way of forcing a
# numeric-looking entry to stay text (sets quotePrefix instead of
# reformatting the cell), so do that for the numeric-looking values.
$ws.Range("A2").Value = "'2"
$ws.Range("B2").Value = "'1"
$ws.Range("A3").Value = "'3"
$ws.Range("B3").Value = "'4"
$ws.Range("A4").Value = "'4"
$ws.Range("B4").Value = "'3"
$ws.Range("A5").Value = "'5"
$ws.Range("B5").Value = "'7"
$ws.Range("A6").Value = "הדס"
$ws.Range("B6").Value = "'2"

# Rows 7-11 are removed entirely (dimension shrinks to A1:B6).
$ws.Range("A7:B11").ClearContents()
